$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B7"  = 5.093999999999999
    "A9"  = -21.93
    "B12" = 5.942399999999997
    "E15" = 16.4664
    "A18" = -22.31950000000002
    "A20" = -19.41299999999999
    "B26" = 4.259000000000005
    "A27" = -21.86959999999999
    "B27" = 5.470200000000003
    "B29" = 5.061599999999997
    "B37" = 8.690900000000006
    "B38" = 4.973999999999999
    "E38" = 16.18939999999999
    "E44" = 16.73889999999999
    "B51" = 5.630000000000003
    "E51" = 17.37460000000001
    "B55" = 5.307299999999999
    "E57" = 16.6801
    "E63" = 18.52510000000001
    "A69" = -21.68189999999998
    "B69" = 5.742099999999994
    "B70" = 6.445900000000002
    "E70" = 17.67700000000001
    "A76" = -19.55049999999999
    "A82" = -22.09120000000001
    "B83" = 5.631
    "E99" = 16.42789999999999
    "B102" = 8.390300000000007
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
